$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6627.25
$ws.Range("J51").Value = 6627.25
$ws.Range("L51").Value = 6627.25
$ws.Range("N51").Value = -7595.25

$ws.Range("H62").Value = 2854.1538
$ws.Range("I62").Value = 2472.1428
$ws.Range("J62").Value = 3299.8333
$ws.Range("K62").Value = 2472.1428
$ws.Range("L62").Value = 3299.8333
$ws.Range("M62").Value = -1848.1428
$ws.Range("N62").Value = -4547.8333

$ws.Range("H65").Value = 2854.1538
$ws.Range("I65").Value = 2472.1428
$ws.Range("J65").Value = 3299.8333
$ws.Range("K65").Value = 12360.714
$ws.Range("L65").Value = 16499.1665
$ws.Range("M65").Value = -9240.714
$ws.Range("N65").Value = -22739.1665

$ws.Range("H98").Value = 902.0769
$ws.Range("I98").Value = 864.9
$ws.Range("K98").Value = 864.9
$ws.Range("M98").Value = 633.1

$ws.Range("H121").Value = 9007.666999999999
$ws.Range("J121").Value = 9593.929
$ws.Range("L121").Value = 28781.787
$ws.Range("N121").Value = -32275.787

$ws.Range("H122").Value = 902.0769
$ws.Range("I122").Value = 864.9
$ws.Range("K122").Value = 2594.7
$ws.Range("M122").Value = -144.6999999999998

$ws.Range("H137").Value = 1799.0278
$ws.Range("I137").Value = 1666.5555
$ws.Range("J137").Value = 2196.4443
$ws.Range("K137").Value = 4999.666499999999
$ws.Range("L137").Value = 6589.3329
$ws.Range("M137").Value = -2449.666499999999
$ws.Range("N137").Value = -11689.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3049.3696
$ws.Range("I32").Value = 2440
$ws.Range("K32").Value = 2440
$ws.Range("M32").Value = -2153

$ws.Range("H61").Value = 785523.5
$ws.Range("I61").Value = 1202115.2
$ws.Range("J61").Value = 4413.875
$ws.Range("K61").Value = 1202115.2
$ws.Range("L61").Value = 4413.875
$ws.Range("M61").Value = -1201903.2
$ws.Range("N61").Value = -4837.875

$ws.Range("H97").Value = 2119.2856
$ws.Range("I97").Value = 1799.091
$ws.Range("J97").Value = 3293.3333
$ws.Range("K97").Value = 1799.091
$ws.Range("L97").Value = 3293.3333
$ws.Range("M97").Value = -1303.091
$ws.Range("N97").Value = -4285.3333

$ws.Range("H122").Value = 1672.2858
$ws.Range("I122").Value = 1664.7273
$ws.Range("J122").Value = 1700
$ws.Range("K122").Value = 4994.1819
$ws.Range("L122").Value = 5100
$ws.Range("M122").Value = -2544.1819
$ws.Range("N122").Value = -10000

$ws.Range("H132").Value = 21529.309
$ws.Range("I132").Value = 2465.6667
$ws.Range("J132").Value = 37869.57
$ws.Range("K132").Value = 7397.000100000001
$ws.Range("L132").Value = 113608.71
$ws.Range("M132").Value = -4867.000100000001
$ws.Range("N132").Value = -118668.71

$ws.Range("H136").Value = 785523.5
$ws.Range("I136").Value = 1202115.2
$ws.Range("J136").Value = 4413.875
$ws.Range("K136").Value = 3606345.6
$ws.Range("L136").Value = 13241.625
$ws.Range("M136").Value = -3603795.6
$ws.Range("N136").Value = -18341.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3733.1924
$ws.Range("I94").Value = 1883.9375
$ws.Range("J94").Value = 6692
$ws.Range("K94").Value = 1883.9375
$ws.Range("L94").Value = 6692
$ws.Range("M94").Value = -1432.9375
$ws.Range("N94").Value = -7594

$ws.Range("H99").Value = 2144.9333
$ws.Range("I99").Value = 1942.4286
$ws.Range("K99").Value = 1942.4286
$ws.Range("M99").Value = -444.4286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9048.42
$ws.Range("I31").Value = 10184.829
$ws.Range("K31").Value = 10184.829
$ws.Range("M31").Value = -9889.829

$ws.Range("H34").Value = 9048.42
$ws.Range("I34").Value = 10184.829
$ws.Range("K34").Value = 10184.829
$ws.Range("M34").Value = -9982.829

$ws.Range("H58").Value = 19714.482
$ws.Range("I58").Value = 1387.9375
$ws.Range("J58").Value = 46371.273
$ws.Range("K58").Value = 1387.9375
$ws.Range("L58").Value = 46371.273
$ws.Range("M58").Value = -1184.9375
$ws.Range("N58").Value = -46777.273

$ws.Range("H136").Value = 19714.482
$ws.Range("I136").Value = 1387.9375
$ws.Range("J136").Value = 46371.273
$ws.Range("K136").Value = 4163.8125
$ws.Range("L136").Value = 139113.819
$ws.Range("M136").Value = -1613.8125
$ws.Range("N136").Value = -144213.819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4000
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H33").Value = 249.5

$ws.Range("H80").Value = 1543.2858
$ws.Range("J80").Value = 1543.2858
$ws.Range("L80").Value = 4629.857400000001
$ws.Range("N80").Value = -6501.857400000001

$ws.Range("H83").Value = 1543.2858
$ws.Range("J83").Value = 1543.2858
$ws.Range("L83").Value = 13889.5722
$ws.Range("N83").Value = -23249.5722

$ws.Range("H96").Value = 5097.8
$ws.Range("J96").Value = 5097.8
$ws.Range("L96").Value = 15293.4
$ws.Range("N96").Value = -19411.4

$ws.Range("H113").Value = 546.8889
$ws.Range("I113").Value = 564.5
$ws.Range("J113").Value = 541.8570999999999
$ws.Range("K113").Value = 1693.5
$ws.Range("L113").Value = 1625.5713
$ws.Range("M113").Value = 476.5
$ws.Range("N113").Value = -5965.5713

$ws.Range("H131").Value = 701.05
$ws.Range("J131").Value = 717.44794
$ws.Range("L131").Value = 2152.34382
$ws.Range("N131").Value = -12232.34382

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1529.2069
$ws.Range("I102").Value = 1534.381
$ws.Range("K102").Value = 1534.381
$ws.Range("M102").Value = 87.61899999999991

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4412.6665
$ws.Range("I40").Value = 3744.5454
$ws.Range("K40").Value = 3744.5454
$ws.Range("M40").Value = -3608.5454

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3247944.2
$ws.Range("J107").Value = 7578090.5
$ws.Range("L107").Value = 22734271.5
$ws.Range("N107").Value = -22738111.5

$ws.Range("H122").Value = 1912.375
$ws.Range("I122").Value = 1725.2858
$ws.Range("K122").Value = 5175.857400000001
$ws.Range("M122").Value = -2725.857400000001

$ws.Range("H132").Value = 3061.4546
$ws.Range("I132").Value = 2772.25
$ws.Range("J132").Value = 3832.6667
$ws.Range("K132").Value = 8316.75
$ws.Range("L132").Value = 11498.0001
$ws.Range("M132").Value = -5786.75
$ws.Range("N132").Value = -16558.0001

$ws.Range("H136").Value = 1326.7632
$ws.Range("I136").Value = 929.75
$ws.Range("K136").Value = 2789.25
$ws.Range("M136").Value = -239.25
